# Updates the cryptos list (Coin / Link / Price / Volume(1h)) to match the
# "Wed Feb 28 15:15:13 UTC 2024" GitHub Actions refresh.
#
# Column D ("Price") cells are stored as literal text in the workbook even
# though most of the values look numeric (e.g. "39.58"). Excel's normal
# Range.Value assignment auto-coerces such strings into real numbers (and
# also tends to stamp a new quote-prefix / text number-format style on the
# cell), which would not match the original plain-text cells. To keep the
# cell as plain text with no style changes, we round-trip the value through
# a helper cell as a quoted-text formula, then Copy / PasteSpecial
# (values-only) it into the destination - that yields a literal string cell
# identical in shape to the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = "H1"

function Set-TextCell($cellRef, $text) {
    $ws.Range($scratch).Value = "=""" + $text + """"
    $ws.Range($scratch).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# Simple rows: only Price (D) and/or Volume(1h) (E) changed.
$rows = @(
    @{ r = 2;  d = "60.873.64"; e = "  +6.92%  " },
    @{ r = 3;  d = "3.349.03";  e = "  +2.95%  " },
    @{ r = 4;  d = "1.00";      e = "  -0.08%  " },
    @{ r = 5;  d = "413.25";    e = "  +4.76%  " },
    @{ r = 6;  d = "111.26";    e = "  +2.32%  " },
    @{ r = 7;  d = "0.587";     e = "  +4.17%  " },
    @{ r = 8;  d = "0.999";     e = "  -0.07%  " },
    @{ r = 9;  d = "0.639";     e = "  +2.72%  " },
    @{ r = 10; d = "39.58";     e = "  +1.12%  " },
    @{ r = 11; d = "0.0991";    e = "  +2.98%  " },
    @{ r = 12; d = $null;       e = "  +1.21%  " },
    @{ r = 13; d = "3.876.31";  e = "  +2.82%  " },
    @{ r = 16; d = "3.323.29";  e = "  +1.85%  " },
    @{ r = 17; d = "1.04";      e = "  +1.37%  " },
    @{ r = 18; d = "60.592.74"; e = "  +6.66%  " },
    @{ r = 19; d = "10.66";     e = "  +0.02%  " },
    @{ r = 20; d = "3.37";      e = "  +2.44%  " },
    @{ r = 21; d = "0.0000110"; e = "  +5.11%  " },
    @{ r = 22; d = "13.09";     e = "  +0.75%  " },
    @{ r = 23; d = "303.14";    e = "  +1.36%  " },
    @{ r = 24; d = "75.04";     e = "  +1.17%  " },
    @{ r = 25; d = "3.20";      e = "  +1.92%  " },
    @{ r = 26; d = "7.92";      e = "  +9.83%  " },
    @{ r = 27; d = "28.73";     e = "  +2.78%  " },
    @{ r = 28; d = "4.49";      e = "  +2.14%  " },
    @{ r = 29; d = "7.98";      e = "  +1.63%  " },
    @{ r = 30; d = $null;       e = "  +6.07%  " },
    @{ r = 33; d = "11.46";     e = "  +4.46%  " },
    @{ r = 34; d = "0.999";     e = "  -0.08%  " },
    @{ r = 35; d = "39.32";     e = "  +3.72%  " },
    @{ r = 36; d = "0.0507";    e = "  +5.43%  " },
    @{ r = 37; d = "52.29";     e = "  +1.60%  " },
    @{ r = 38; d = "3.13";      e = "  +1.63%  " },
    @{ r = 39; d = "0.999";     e = "  -0.16%  " },
    @{ r = 40; d = "3.41";      e = "  -1.78%  " },
    @{ r = 43; d = "0.123";     e = "  +3.28%  " },
    @{ r = 44; d = "1.91";      e = "  +0.08%  " },
    @{ r = 45; d = "3.95";      e = "  -2.05%  " },
    @{ r = 46; d = "16.91";     e = "  -1.27%  " },
    @{ r = 47; d = $null;       e = "  +8.56%  " },
    @{ r = 48; d = "22.47";     e = "  +2.58%  " },
    @{ r = 49; d = "2.180.36";  e = "  +1.87%  " },
    @{ r = 50; d = $null;       e = "  +1.92%  " },
    @{ r = 51; d = "1.99";      e = "  -1.66%  " }
)

foreach ($row in $rows) {
    if ($null -ne $row.d) {
        Set-TextCell ("D" + $row.r) $row.d
    }
    $ws.Range("E" + $row.r).Value = $row.e
}

# Rows that reshuffled (coin order changed) - Coin / Link / Price / Volume(1h)
# all move together.

# Row 14/15: Polkadot <-> Chainlink swapped order
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D14" "20.00"
$ws.Range("E14").Value = "  +5.12%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D15" "8.40"
$ws.Range("E15").Value = "  +2.90%  "

# Row 31/32: Toncoin <-> Hedera swapped order
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D31" "0.116"
$ws.Range("E31").Value = "  +5.21%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D32" "2.63"
$ws.Range("E32").Value = "  +24.25%  "

# Row 41/42: Monero <-> TheGraph swapped order
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D41" "0.301"
$ws.Range("E41").Value = "  +6.57%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D42" "137.66"
$ws.Range("E42").Value = "  +2.49%  "

# Clean up the scratch cell used for the text round-trip.
$ws.Range($scratch).Clear()
